$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1877133105802048
$ws.Range("C2").Value = 0.5767918088737202
$ws.Range("J2").Value = 0.0136518771331058
$ws.Range("P2").Value = 0.1160409556313993
$ws.Range("S2").Value = 0.10580204778157
$ws.Range("B3").Value = 0.005882352941176471
$ws.Range("C3").Value = 0.02352941176470588
$ws.Range("J3").Value = 0.02941176470588235
$ws.Range("P3").Value = 0.7352941176470589
$ws.Range("S3").Value = 0.2058823529411765
$ws.Range("J4").Value = 0.1063829787234043
$ws.Range("P4").Value = 0.6170212765957447
$ws.Range("S4").Value = 0.2765957446808511
$ws.Range("B6").Value = 0.05681818181818182
$ws.Range("D6").Value = 0.01515151515151515
$ws.Range("F6").Value = 0.05303030303030303
$ws.Range("J6").Value = 0.3106060606060606
$ws.Range("O6").Value = 0.01515151515151515
$ws.Range("Q6").Value = 0.1553030303030303
$ws.Range("R6").Value = 0.07575757575757576
$ws.Range("S6").Value = 0.3181818181818182
$ws.Range("B7").Value = 0.06914893617021277
$ws.Range("D7").Value = 0.02659574468085106
$ws.Range("F7").Value = 0.04787234042553191
$ws.Range("J7").Value = 0.1170212765957447
$ws.Range("O7").Value = 0.02127659574468085
$ws.Range("Q7").Value = 0.2180851063829787
$ws.Range("R7").Value = 0.09042553191489362
$ws.Range("S7").Value = 0.4095744680851064
$ws.Range("B8").Value = 0.09579439252336448
$ws.Range("D8").Value = 0.02336448598130841
$ws.Range("F8").Value = 0.07476635514018691
$ws.Range("J8").Value = 0.08644859813084112
$ws.Range("O8").Value = 0.01635514018691589
$ws.Range("Q8").Value = 0.1869158878504673
$ws.Range("R8").Value = 0.1121495327102804
$ws.Range("S8").Value = 0.4042056074766355
$ws.Range("B9").Value = 0.0821917808219178
$ws.Range("D9").Value = 0.0091324200913242
$ws.Range("F9").Value = 0.0502283105022831
$ws.Range("J9").Value = 0.0684931506849315
$ws.Range("O9").Value = 0.0273972602739726
$ws.Range("Q9").Value = 0.2191780821917808
$ws.Range("R9").Value = 0.0821917808219178
$ws.Range("S9").Value = 0.4611872146118721
$ws.Range("B10").Value = 0.1185243328100471
$ws.Range("D10").Value = 0.02197802197802198
$ws.Range("F10").Value = 0.07378335949764521
$ws.Range("J10").Value = 0.1004709576138148
$ws.Range("O10").Value = 0.01648351648351648
$ws.Range("Q10").Value = 0.2182103610675039
$ws.Range("R10").Value = 0.08869701726844584
$ws.Range("S10").Value = 0.3618524332810047
$ws.Range("G11").Value = 0.1597222222222222
$ws.Range("J11").Value = 0.08333333333333333
$ws.Range("K11").Value = 0.2118055555555556
$ws.Range("L11").Value = 0.5173611111111112
$ws.Range("S11").Value = 0.02777777777777778
$ws.Range("G12").Value = 0.740506329113924
$ws.Range("J12").Value = 0.2088607594936709
$ws.Range("L12").Value = 0.0379746835443038
$ws.Range("S12").Value = 0.01265822784810127
$ws.Range("F13").Value = 0.02222222222222222
$ws.Range("G13").Value = 0.6666666666666666
$ws.Range("J13").Value = 0.2888888888888889
$ws.Range("S13").Value = 0.02222222222222222
$ws.Range("F15").Value = 0.03474903474903475
$ws.Range("H15").Value = 0.1467181467181467
$ws.Range("I15").Value = 0.06563706563706563
$ws.Range("J15").Value = 0.3783783783783784
$ws.Range("K15").Value = 0.05405405405405406
$ws.Range("M15").Value = 0.01158301158301158
$ws.Range("O15").Value = 0.06563706563706563
$ws.Range("S15").Value = 0.2432432432432433
$ws.Range("F16").Value = 0.03314917127071823
$ws.Range("H16").Value = 0.1104972375690608
$ws.Range("I16").Value = 0.09392265193370165
$ws.Range("J16").Value = 0.3922651933701657
$ws.Range("K16").Value = 0.09944751381215469
$ws.Range("M16").Value = 0.02209944751381215
$ws.Range("O16").Value = 0.09944751381215469
$ws.Range("S16").Value = 0.1491712707182321
$ws.Range("F17").Value = 0.03112033195020747
$ws.Range("H17").Value = 0.1701244813278008
$ws.Range("I17").Value = 0.09336099585062241
$ws.Range("J17").Value = 0.3858921161825726
$ws.Range("K17").Value = 0.07468879668049792
$ws.Range("M17").Value = 0.01867219917012448
$ws.Range("O17").Value = 0.07261410788381743
$ws.Range("S17").Value = 0.1535269709543569
$ws.Range("F18").Value = 0.04205607476635514
$ws.Range("H18").Value = 0.1308411214953271
$ws.Range("I18").Value = 0.09813084112149532
$ws.Range("J18").Value = 0.4252336448598131
$ws.Range("K18").Value = 0.08878504672897196
$ws.Range("M18").Value = 0.009345794392523364
$ws.Range("O18").Value = 0.07476635514018691
$ws.Range("S18").Value = 0.1308411214953271
$ws.Range("F19").Value = 0.02304832713754647
$ws.Range("H19").Value = 0.1955390334572491
$ws.Range("I19").Value = 0.08847583643122676
$ws.Range("J19").Value = 0.3561338289962825
$ws.Range("K19").Value = 0.1018587360594796
$ws.Range("M19").Value = 0.02156133828996282
$ws.Range("N19").Value = 0.001486988847583643
$ws.Range("O19").Value = 0.07657992565055761
$ws.Range("S19").Value = 0.1353159851301115
